{"js": "// The rendered \"Date\" paragraph (e.g. \"2017-08-16 16:49:34\") holds the\n// render time-stamp text in its own trailing run. The R Markdown -> docx\n// rendering pipeline was re-run (per the commit message, fixing the\n// html_notebook / full-site rendering), producing a new render\n// time-stamp. Update the visible time-stamp text in place, leaving the\n// surrounding runs (and everything else in the document) untouched.\nconst oldTime = \"16:49:34\";\nconst newTime = \"18:06:47\";\n\nconst results = context.document.body.search(oldTime, {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  for (const range of results.items) {\n    range.insertText(newTime, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The rendered \"Date\" line (e.g. \"2017-08-16 16:49:34\") stores the\n# render time-stamp text in its own run. The R Markdown -> docx\n# pipeline was re-run (per the commit message, fixing html_notebook /\n# full-site rendering), producing a new render time-stamp. Update the\n# visible time-stamp text in place, leaving everything else (including\n# the neighboring runs) untouched.\n\n$d = $word.ActiveDocument\n\n$oldTime = \"16:49:34\"\n$newTime = \"18:06:47\"\n\n$target = $d.Content\n$find = $target.Find\n$find.ClearFormatting()\n$find.Text = $oldTime\n$find.Forward = $true\n$find.Wrap = 0\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute()\nif ($found) {\n    $target.Delete()\n    $target.InsertAfter($newTime)\n}\n"}
